# PlayerPerformance_6666.xlsx update script
# - Adds a new "Player Info" worksheet as the first sheet
# - Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling" sheets
# - Replaces the full match-card URL values with the bare match code (4603)

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- Create the new "Player Info" worksheet, placed before "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row (keep ID as text, matching the source data which stores every
# value - even numeric-looking ones - as text)
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6666"
$playerInfo.Range("B2").Value = "Pramod Madushan Liyanagamage"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# Match header styling used elsewhere in the workbook (bold font, thin border,
# centered horizontally, top-aligned vertically)
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$playerInfo.Activate()
$playerInfo.Range("A1").Select()

# Re-fetch sheet references by name, since adding a worksheet can invalidate
# previously captured references.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4603"

# --- Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4603"
